$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 814, shifting existing rows 814:855 down to 815:856.
$ws.Rows.Item(814).Insert()

# Copy the date/weekday text from the row above (same "2026/02/15" / "日"
# entry group) so the new cells stay plain text instead of being
# auto-parsed as a date when typed directly.
$ws.Range("A813:B813").Copy()
$ws.Range("A814").PasteSpecial()

# Fill in the new log entry's time and ranking values.
$ws.Cells.Item(814, 3).Value = 13
$ws.Cells.Item(814, 4).Value = 201
